# SCHRÄG BOM — "Added back the FM normalization. This time with a bit more
# than 1 OCT range." (#8)
#
# A new SMD resistor line (R11, 2M2) is inserted into the BOM as the new
# row 27, pushing every row below it (U1/U2/U4/U5, U3, the "THT Parts"
# subtotal row, and all THT component rows) down by one. Row/merge refs
# and the two CONCAT/SUM subtotal formulas are re-pointed automatically by
# the row-insert. Separately, the existing 100K resistor group (row 23)
# had its quantity corrected from 15 to 14 (14 designators are listed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new BOM row for R11 right above the former row 27 (U1,U2,U4,U5) ---
$ws.Rows.Item(27).Insert()

$ws.Range("A27").Value = "R11"
$ws.Range("B27").Value = "Resistor, 1%"
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = "2M2"
$ws.Range("G27").Value = "Panasonic"
$ws.Range("H27").Value = "ERJ-3EKF2204V"
$ws.Range("I27").Value = "667-ERJ-3EKF2204V"

# "0603" round-trips through Range.Value as the number 603 (loses the
# leading zero), so clone the existing "0603" text cell (C4 - same package
# used by every other 0603 resistor/capacitor row) instead of assigning it
# directly. This also picks up the matching "s=4" style in one shot.
$ws.Range("C4").Copy()
$ws.Range("C27").PasteSpecial()

# Match the style used by the other resistor rows (e.g. row 14/26): bold-ish
# "s=4" style on B/G (C27 already got its style from the paste above).
$ws.Range("B27").Style = $ws.Range("B26").Style
$ws.Range("G27").Style = $ws.Range("G26").Style

# --- Fix the pre-existing quantity bug on the 100K resistor group (still row 23) ---
$ws.Range("E23").Value = 14
